$wb = $excel.ActiveWorkbook

# --- Rename the second worksheet tab ---
$ws2 = $wb.Worksheets.Item("Include from Hearing Observat")
$ws2.Name = "Include from LOINC"

# --- Update Metadata sheet (sheet1): Date and Description values ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B8").Value = "2022-03-24T23:37:25-04:00"
$ws1.Range("B13").Value = "SPLASCH hearing observation profile: codes representing hearing questions that are answered by the ValueSet SPLASCHHearingObservationValueCodeableConceptVS"

# --- "Include from LOINC" sheet: rebuild the Concept table ---

# Give the new header cell B1 the same (bold) formatting as A1
$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("B1").PasteSpecial(-4122) | Out-Null

# B2 doesn't exist yet either - format it like A2/B3
$ws2.Range("B3").Copy() | Out-Null
$ws2.Range("B2").PasteSpecial(-4122) | Out-Null

# Extend the body with new formatted rows 5-7, copied from row 4's format
$ws2.Range("A4:B4").Copy() | Out-Null
$ws2.Range("A5:B7").PasteSpecial(-4122) | Out-Null

# Header row
$ws2.Range("A1").Value = "Concept"
$ws2.Range("B1").Value = "Description"

# Row 2
$ws2.Range("A2").Value = "95744-9"
$ws2.Range("B2").Value = "Hearing.ability to hear during assessment period [CMS Assessment]"

# Row 3
$ws2.Range("A3").Value = "54599-6"
$ws2.Range("B3").Value = "Hearing aid used during assessment period [CMS Assessment]"

# Row 4
$ws2.Range("A4").Value = "67235-2"
$ws2.Range("B4").Value = "Does the examinee have hearing aids that cannot be removed [PhenX]"

# Row 5
$ws2.Range("A5").Value = "94900-8"
$ws2.Range("B5").Value = "Need for and availability of a hearing aid [CMS Assessment]"

# Row 6 is intentionally blank (keeps formatting, no text)
$ws2.Range("A6").Value = ""
$ws2.Range("B6").Value = ""

# Row 7
$ws2.Range("A7").Value = "System URI"
$ws2.Range("B7").Value = "http://loinc.org"
